$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing row label in A8 from "descuento_aplicado_pct" to "descuento_pct"
$ws.Range("A8").Value = "descuento_pct"

# Add new row 9 data for "monto_final"
$ws.Range("A9").Value = "monto_final"
$ws.Range("B9").Value = 39.06
$ws.Range("C9").Value = 119.04
$ws.Range("D9").Value = 79.98
$ws.Range("E9").Value = -80.91
$ws.Range("F9").Value = 239
$ws.Range("G9").Value = 8
$ws.Range("H9").Value = 0.4
